$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos -- update B10/C10 to new Portuguese objectives text ---
$ws.Range("B10").Value = 'Apresentar ao aluno o conceito de uma organização e os fundamentos de sua administração;Caracterizar as diversas áreas funcionais existentes nas organizações;Despertar o interesse dos alunos para questões de gestão'
$ws.Range("C10").Value = 'Apresentar ao aluno o conceito de uma organização e os fundamentos de sua administração;Caracterizar as diversas áreas funcionais existentes nas organizações;Despertar o interesse dos alunos para questões de gestão'

# --- Row 13: drop A13 label, set B13/C13 to the professor name (moved up from old row 13 B/C) ---
$ws.Range("A13").Clear()
$ws.Range("B13").Value = '849935 - Humberto Felipe da Silva'
$ws.Range("C13").Value = '849935 - Humberto Felipe da Silva'
$ws.Rows.Item(13).EntireRow.AutoFit()

# --- Row 14: Programa resumido -- new short Portuguese syllabus text ---
$ws.Range("B14").Value = '1 - A Administração das Organizações. 2 - O processo administrativo. 3 – Processos de Gestão'
$ws.Range("C14").Value = '1 - A Administração das Organizações. 2 - O processo administrativo. 3 – Processos de Gestão'
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: Short syllabus -- row height 120 -> 60 ---
$ws.Rows.Item(15).RowHeight = 60

# --- Row 16: Programa -- new long Portuguese program text ---
$ws.Range("B16").Value = '1 - A Administração das organizações - definindo a administração2 - O processo administrativo: planejamento, organização, direção, controle3 – Processos de Gestão: Marketing, Finanças, Gestão de Pessoas, Produção e Operações, Pesquisa e Desenvolvimento, Tecnologia da Informação, Logística e Meio Ambiente.'
$ws.Range("C16").Value = '1 - A Administração das organizações - definindo a administração2 - O processo administrativo: planejamento, organização, direção, controle3 – Processos de Gestão: Marketing, Finanças, Gestão de Pessoas, Produção e Operações, Pesquisa e Desenvolvimento, Tecnologia da Informação, Logística e Meio Ambiente.'

# --- Row 17: Syllabus -- add B17/C17 English long syllabus text, set height 120 ---
$ws.Range("B17").Value = '- The Administration of organizations - defining the administration 2 - The administrative process: planning, organization, direction, control 3 - Management Processes: Marketing, Finance, People Management, Production and Operations, Research and Development, Information Technology, Logistics and Environment.'
$ws.Range("C17").Value = '- The Administration of organizations - defining the administration 2 - The administrative process: planning, organization, direction, control 3 - Management Processes: Marketing, Finance, People Management, Production and Operations, Research and Development, Information Technology, Logistics and Environment.'
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(17).RowHeight = 120

# --- Row 18: Avaliacao -- remove B18/C18, remove custom height ---
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).EntireRow.AutoFit()

# --- Row 19: Metodo -- new avaliacao method text ---
$ws.Range("B19").Value = 'O sistema de avaliação será continuo com aplicação de avaliações escritas utilizando-se avaliações em papel como em sistema informacionais, ademais serão realizados seminários, projetos, entrega de trabalho em formato de artigo e Estudos de Casos.'
$ws.Range("C19").Value = 'O sistema de avaliação será continuo com aplicação de avaliações escritas utilizando-se avaliações em papel como em sistema informacionais, ademais serão realizados seminários, projetos, entrega de trabalho em formato de artigo e Estudos de Casos.'

# --- Row 21: Norma de recuperacao -- row height 120 -> 60 ---
$ws.Rows.Item(21).RowHeight = 60

# --- Row 22 (new): Bibliografia -- add new row with bibliography text ---
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'LEMOS, Paulo de Mattos et al. Gestão estratégica de empresas. Rio de Janeiro: Fundação Getúlio Vargas, 2014.Ludovico, Nelson. Gestão estratégica de negócios. São Paulo: Saraiva, 2018Serra, Fernando Ribeiro et al. Gestão estratégica: conceitos e casos. São Paulo: Atlas, 2014.'
$ws.Range("C22").Value = 'LEMOS, Paulo de Mattos et al. Gestão estratégica de empresas. Rio de Janeiro: Fundação Getúlio Vargas, 2014.Ludovico, Nelson. Gestão estratégica de negócios. São Paulo: Saraiva, 2018Serra, Fernando Ribeiro et al. Gestão estratégica: conceitos e casos. São Paulo: Atlas, 2014.'
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(22).RowHeight = 120

# --- Columns: split col A/B merged width definition (A keeps 30.7109375, B keeps its own 60.7109375) ---
$ws.Columns.Item(2).ColumnWidth = 60.7109375
